$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Grupo 1 - Item 1 - teste 1 (1)"
$ws.Range("A2").Value = "teste descricao detalhada 1"
$ws.Range("A5").Value = "Grupo 1 - Item 2 - teste 2 (2)"
$ws.Range("A6").Value = "teste descricao detalhada 2"
$ws.Range("A9").Value = "Grupo 1 - Item 3 - teste 3 (3)"
$ws.Range("A10").Value = "teste descricao detalhada 3"
$ws.Range("A13").Value = "Grupo 1 - Item 4 - teste 4 (4)"
$ws.Range("A14").Value = "teste descricao detalhada 4"
$ws.Range("A17").Value = "Grupo 1 - Item 5 - teste 5 (5)"
$ws.Range("A18").Value = "teste descricao detalhada 5"
$ws.Range("A21").Value = "Grupo 1 - Item 6 - teste 6 (6)"
$ws.Range("A22").Value = "teste descricao detalhada 6"
$ws.Range("A25").Value = "Grupo 1 - Item 7 - teste 7 (7)"
$ws.Range("A26").Value = "teste descricao detalhada 7"
$ws.Range("A29").Value = "Grupo 1 - Item 8 - teste 8 (8)"
$ws.Range("A30").Value = "teste descricao detalhada 8"
$ws.Range("A33").Value = "Grupo 1 - Item 9 - teste 9 (9)"
$ws.Range("A34").Value = "teste descricao detalhada 9"
$ws.Range("A37").Value = "Grupo 1 - Item 10 - teste 10 (10)"
$ws.Range("A38").Value = "teste descricao detalhada 10"
$ws.Range("A41").Value = "Grupo 1 - Item 11 - teste 11 (11)"
$ws.Range("A42").Value = "teste descricao detalhada 11"
$ws.Range("A45").Value = "Grupo 1 - Item 12 - teste 12 (12)"
$ws.Range("A46").Value = "teste descricao detalhada 12"
$ws.Range("A49").Value = "Grupo 1 - Item 13 - teste 13 (13)"
$ws.Range("A50").Value = "teste descricao detalhada 13"
$ws.Range("A53").Value = "Grupo 1 - Item 14 - teste 14 (14)"
$ws.Range("A54").Value = "teste descricao detalhada 14"
$ws.Range("A57").Value = "Grupo 1 - Item 15 - teste 15 (15)"
$ws.Range("A58").Value = "teste descricao detalhada 15"
$ws.Range("A61").Value = "Grupo 1 - Item 16 - teste 16 (16)"
$ws.Range("A62").Value = "teste descricao detalhada 16"
$ws.Range("A65").Value = "Grupo 1 - Item 17 - teste 17 (17)"
$ws.Range("A66").Value = "teste descricao detalhada 17"
$ws.Range("A69").Value = "Grupo 1 - Item 18 - teste 18 (18)"
$ws.Range("A70").Value = "teste descricao detalhada 18"
$ws.Range("A73").Value = "Grupo 1 - Item 19 - teste 19 (19)"
$ws.Range("A74").Value = "teste descricao detalhada 19"
$ws.Range("A77").Value = "Grupo 1 - Item 20 - teste 20 (20)"
$ws.Range("A78").Value = "teste descricao detalhada 20"
$ws.Range("A81").Value = "Grupo 1 - Item 21 - teste 21 (21)"
$ws.Range("A82").Value = "teste descricao detalhada 21"
$ws.Range("A85").Value = "Grupo 1 - Item 22 - teste 22 (22)"
$ws.Range("A86").Value = "teste descricao detalhada 22"
$ws.Range("A89").Value = "Grupo 1 - Item 23 - teste 23 (23)"
$ws.Range("A90").Value = "teste descricao detalhada 23"
$ws.Range("A93").Value = "Grupo 1 - Item 24 - teste 24 (24)"
$ws.Range("A94").Value = "teste descricao detalhada 24"
$ws.Range("A97").Value = "Grupo 1 - Item 25 - teste 25 (25)"
$ws.Range("A98").Value = "teste descricao detalhada 25"
$ws.Range("A101").Value = "Grupo 1 - Item 26 - teste 26 (26)"
$ws.Range("A102").Value = "teste descricao detalhada 26"
$ws.Range("A105").Value = "Grupo 2 - Item 27 - teste 27 (27)"
$ws.Range("A106").Value = "teste descricao detalhada 27"
$ws.Range("A109").Value = "Grupo 2 - Item 28 - teste 28 (28)"
$ws.Range("A110").Value = "teste descricao detalhada 28"
$ws.Range("A113").Value = "Grupo 2 - Item 29 - teste 29 (29)"
$ws.Range("A114").Value = "teste descricao detalhada 29"
$ws.Range("A117").Value = "Grupo 2 - Item 30 - teste 30 (30)"
$ws.Range("A118").Value = "teste descricao detalhada 30"
$ws.Range("A121").Value = "Grupo 2 - Item 31 - teste 31 (31)"
$ws.Range("A122").Value = "teste descricao detalhada 31"
$ws.Range("A125").Value = "Grupo 2 - Item 32 - teste 32 (32)"
$ws.Range("A126").Value = "teste descricao detalhada 32"
$ws.Range("A129").Value = "Grupo 2 - Item 33 - teste 33 (33)"
$ws.Range("A130").Value = "teste descricao detalhada 33"
$ws.Range("A133").Value = "Grupo 2 - Item 34 - teste 34 (34)"
$ws.Range("A134").Value = "teste descricao detalhada 34"
$ws.Range("A137").Value = "Grupo 2 - Item 35 - teste 35 (35)"
$ws.Range("A138").Value = "teste descricao detalhada 35"
$ws.Range("A141").Value = "Grupo 2 - Item 36 - teste 36 (36)"
$ws.Range("A142").Value = "teste descricao detalhada 36"
$ws.Range("A145").Value = "Grupo 2 - Item 37 - teste 37 (37)"
$ws.Range("A146").Value = "teste descricao detalhada 37"
$ws.Range("A149").Value = "Grupo 2 - Item 38 - teste 38 (38)"
$ws.Range("A150").Value = "teste descricao detalhada 38"
$ws.Range("A153").Value = "Grupo 2 - Item 39 - teste 39 (39)"
$ws.Range("A154").Value = "teste descricao detalhada 39"
$ws.Range("A157").Value = "Grupo 2 - Item 40 - teste 40 (40)"
$ws.Range("A158").Value = "teste descricao detalhada 40"
$ws.Range("A161").Value = "Grupo 2 - Item 41 - teste 41 (41)"
$ws.Range("A162").Value = "teste descricao detalhada 41"
$ws.Range("A165").Value = "Grupo 2 - Item 42 - teste 42 (42)"
$ws.Range("A166").Value = "teste descricao detalhada 42"
$ws.Range("A169").Value = "Grupo 2 - Item 43 - teste 43 (43)"
$ws.Range("A170").Value = "teste descricao detalhada 43"
$ws.Range("A173").Value = "Grupo 2 - Item 44 - teste 44 (44)"
$ws.Range("A174").Value = "teste descricao detalhada 44"
$ws.Range("A177").Value = "Grupo 2 - Item 45 - teste 45 (45)"
$ws.Range("A178").Value = "teste descricao detalhada 45"
$ws.Range("A181").Value = "Grupo 2 - Item 46 - teste 46 (46)"
$ws.Range("A182").Value = "teste descricao detalhada 46"
$ws.Range("A185").Value = "Grupo 2 - Item 47 - teste 47 (47)"
$ws.Range("A186").Value = "teste descricao detalhada 47"
$ws.Range("A189").Value = "Grupo 2 - Item 48 - teste 48 (48)"
$ws.Range("A190").Value = "teste descricao detalhada 48"
$ws.Range("A193").Value = "Grupo 2 - Item 49 - teste 49 (49)"
$ws.Range("A194").Value = "teste descricao detalhada 49"
$ws.Range("A197").Value = "Grupo 2 - Item 50 - teste 50 (50)"
$ws.Range("A198").Value = "teste descricao detalhada 50"
$ws.Range("A201").Value = "Grupo 2 - Item 51 - teste 51 (51)"
$ws.Range("A202").Value = "teste descricao detalhada 51"
$ws.Range("A205").Value = "Grupo 2 - Item 52 - teste 52 (52)"
$ws.Range("A206").Value = "teste descricao detalhada 52"
